$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update E2 value from "12345" to "1234567"
$ws.Range("E2").Value = '"1234567"'

# Update the active selection to E2 (as reflected in the saved view state)
$ws.Range("E2").Select()
